$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update last-updated dates (as Excel serial date values) for LLY (Biopharma), AMD (Hardware & Semis), ABBV (Software)
$ws.Range("D5").Value = 45048   # LLY - Biopharma updated to 5/2/2023
$ws.Range("D10").Value = 45048  # AMD - Hardware & Semis updated to 5/2/2023
$ws.Range("D11").Value = 45047  # ABBV - Software updated to 5/1/2023

# Move active selection to D6
$ws.Activate()
$ws.Range("D6").Select()
